# Applies the AFG PlayerPerformance_5936 edit:
#  - Adds two new worksheets: "ODI Batting Extra" and "ODI Bowling Extra"
#    with match-level aggregate stats.
#  - Clears the (empty) INNING_NUMBER cells B5:B7 on "ODI Batting" so the
#    cells are removed entirely rather than present-but-blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clear stray empty cells B5, B6, B7 on the "ODI Batting" sheet so
#    they are fully removed (no <c> element at all), matching the diff.
# ---------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Cells.Item(5, 2).Value = ""
$odiBatting.Cells.Item(6, 2).Value = ""
$odiBatting.Cells.Item(7, 2).Value = ""

# ---------------------------------------------------------------------
# Helper style reference: header rows on the existing sheets use bold
# text, a thin box border all around, centered horizontally and
# top-aligned vertically. Recreate that same look on the header rows
# of the two new sheets.
# ---------------------------------------------------------------------
function Format-HeaderRow($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# Helper: write a value that must stay a text string even though it
# looks numeric / percent-like (Excel would otherwise auto-convert it
# to a number). Forcing NumberFormat to Text ("@") before assignment
# keeps it as a string cell, matching t="inlineStr" cells in the diff.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------
# 2) Add "ODI Batting Extra" sheet (after the last existing sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

# Header row
Set-TextValue $battingExtra.Cells.Item(1, 1) "MATCH_CODE"
Set-TextValue $battingExtra.Cells.Item(1, 2) "BATTING_POSITION"
Set-TextValue $battingExtra.Cells.Item(1, 3) "NUM_4"
Set-TextValue $battingExtra.Cells.Item(1, 4) "NUM_6"
Set-TextValue $battingExtra.Cells.Item(1, 5) "PERCENT_RUNS_OF_TOTAL"
Set-TextValue $battingExtra.Cells.Item(1, 6) "MAN_OF_MATCH"
Format-HeaderRow $battingExtra.Range("A1:F1")

# Row 2 - match 4444
Set-TextValue $battingExtra.Cells.Item(2, 1) "4444"
$battingExtra.Cells.Item(2, 2).Value = 8
Set-TextValue $battingExtra.Cells.Item(2, 3) "0"
Set-TextValue $battingExtra.Cells.Item(2, 4) "0"
Set-TextValue $battingExtra.Cells.Item(2, 5) "0.70%"
Set-TextValue $battingExtra.Cells.Item(2, 6) "NO"

# Row 3 - match 4530
Set-TextValue $battingExtra.Cells.Item(3, 1) "4530"
$battingExtra.Cells.Item(3, 2).Value = 6
Set-TextValue $battingExtra.Cells.Item(3, 3) "2"
Set-TextValue $battingExtra.Cells.Item(3, 4) "0"
Set-TextValue $battingExtra.Cells.Item(3, 5) "5.91%"
Set-TextValue $battingExtra.Cells.Item(3, 6) "NO"

# Row 4 - match 4538 (only MATCH_CODE + MAN_OF_MATCH populated)
Set-TextValue $battingExtra.Cells.Item(4, 1) "4538"
Set-TextValue $battingExtra.Cells.Item(4, 6) "NO"

# Row 5 - match 4539
Set-TextValue $battingExtra.Cells.Item(5, 1) "4539"
$battingExtra.Cells.Item(5, 2).Value = 10
Set-TextValue $battingExtra.Cells.Item(5, 6) "NO"

# Row 6 - match 4582
Set-TextValue $battingExtra.Cells.Item(6, 1) "4582"
$battingExtra.Cells.Item(6, 2).Value = 8
Set-TextValue $battingExtra.Cells.Item(6, 6) "NO"

# Row 7 - match 4585
Set-TextValue $battingExtra.Cells.Item(7, 1) "4585"
$battingExtra.Cells.Item(7, 2).Value = 7
Set-TextValue $battingExtra.Cells.Item(7, 6) "NO"

# Row 8 - match 4588
Set-TextValue $battingExtra.Cells.Item(8, 1) "4588"
$battingExtra.Cells.Item(8, 2).Value = 7
Set-TextValue $battingExtra.Cells.Item(8, 3) "2"
Set-TextValue $battingExtra.Cells.Item(8, 4) "0"
Set-TextValue $battingExtra.Cells.Item(8, 5) "6.57%"
Set-TextValue $battingExtra.Cells.Item(8, 6) "NO"

$battingExtra.Range("A1").Select()

# ---------------------------------------------------------------------
# 3) Add "ODI Bowling Extra" sheet (after "ODI Batting Extra").
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet2)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row
Set-TextValue $bowlingExtra.Cells.Item(1, 1) "MATCH_CODE"
Set-TextValue $bowlingExtra.Cells.Item(1, 2) "MAIDEN_OVERS"
Set-TextValue $bowlingExtra.Cells.Item(1, 3) "PERCENT_WICKETS_OF_ALL"
Format-HeaderRow $bowlingExtra.Range("A1:C1")

# Row 2 - match 4444 (PERCENT_WICKETS_OF_ALL present but blank)
Set-TextValue $bowlingExtra.Cells.Item(2, 1) "4444"
Set-TextValue $bowlingExtra.Cells.Item(2, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(2, 3) ""

# Row 3 - match 4530
Set-TextValue $bowlingExtra.Cells.Item(3, 1) "4530"
Set-TextValue $bowlingExtra.Cells.Item(3, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(3, 3) "10.00%"

# Row 4 - match 4538 (MAIDEN_OVERS + PERCENT present but blank)
Set-TextValue $bowlingExtra.Cells.Item(4, 1) "4538"
Set-TextValue $bowlingExtra.Cells.Item(4, 2) ""
Set-TextValue $bowlingExtra.Cells.Item(4, 3) ""

# Row 5 - match 4539
Set-TextValue $bowlingExtra.Cells.Item(5, 1) "4539"
Set-TextValue $bowlingExtra.Cells.Item(5, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(5, 3) "10.00%"

# Row 6 - match 4582
Set-TextValue $bowlingExtra.Cells.Item(6, 1) "4582"
Set-TextValue $bowlingExtra.Cells.Item(6, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(6, 3) "10.00%"

# Row 7 - match 4585 (PERCENT present but blank)
Set-TextValue $bowlingExtra.Cells.Item(7, 1) "4585"
Set-TextValue $bowlingExtra.Cells.Item(7, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(7, 3) ""

# Row 8 - match 4588
Set-TextValue $bowlingExtra.Cells.Item(8, 1) "4588"
Set-TextValue $bowlingExtra.Cells.Item(8, 2) "0"
Set-TextValue $bowlingExtra.Cells.Item(8, 3) "10.00%"

$bowlingExtra.Range("A1").Select()

# Restore the originally active sheet/selection (first sheet, cell A1)
# so the workbook view stays the same as before the edit.
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()
$playerInfo.Range("A1").Select()
